# feat: add 2022-Q1 data
#
# Before: sheets are 2020-Q4, 2021-Q2, 2021-Q3, 2021-Q4, 总计 (总计 = sheetId 5 / rId5).
# After:  sheets are 2020-Q4, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#   - "2022-Q1" reuses the old "总计" sheet slot (sheetId 5 / rId5) and is repopulated
#     with the new fund-holdings table.
#   - A brand-new "总计" sheet (sheetId 6 / rId6) is inserted after it, holding the same
#     rollup table as before plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1) Build the new "总计" sheet (placed right after the old one) by duplicating
#    the old 总计 layout/formatting, then splice in the 2022-Q1 row at the top.
# ---------------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add($null, $oldTotal)

$newTotal.PageSetup.LeftMargin   = $oldTotal.PageSetup.LeftMargin
$newTotal.PageSetup.RightMargin  = $oldTotal.PageSetup.RightMargin
$newTotal.PageSetup.TopMargin    = $oldTotal.PageSetup.TopMargin
$newTotal.PageSetup.BottomMargin = $oldTotal.PageSetup.BottomMargin
$newTotal.PageSetup.HeaderMargin = $oldTotal.PageSetup.HeaderMargin
$newTotal.PageSetup.FooterMargin = $oldTotal.PageSetup.FooterMargin

# header row (B1:D1 -- column A has no header cell in this table)
$oldTotal.Range("B1:D1").Copy($newTotal.Range("B1:D1"))

# old data rows (2-5) shift down to rows 3-6, carrying their formatting along
$oldTotal.Range("A2:D5").Copy($newTotal.Range("A3:D6"))

# new leading row: reuse A2's formatting (index style) from the old sheet
$oldTotal.Range("A2").Copy($newTotal.Range("A2"))
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 18
$newTotal.Range("D2").Value = 5.91

# re-sequence the row-index column for the rows that moved down
for ($r = 3; $r -le 6; $r++) {
    $newTotal.Range("A$r").Value = $r - 2
}

$newTotal.Name = "__new_总计__"

# ---------------------------------------------------------------------------
# 2) Turn the old "总计" sheet into "2022-Q1" and replace its contents with the
#    fund-holdings table.
# ---------------------------------------------------------------------------
$oldTotal.Name = "2022-Q1"
$ws = $oldTotal

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"
$ws.Range("B1:H1").Font.Bold = $true
$ws.Range("B1:H1").Borders.LineStyle = 1
$ws.Range("B1:H1").HorizontalAlignment = -4108
$ws.Range("B1:H1").VerticalAlignment = -4160

$rows = @(
    @("002168", "嘉实智能汽车股票", "52.15", "92.00", "3.95", "2.0599", 10),
    @("009683", "汇添富创新增长一年定期开放混合A", "18.11", "71.21", "3.01", "0.5451", 7),
    @("009697", "华夏成长精选6个月定期开放混合A", "11.44", "88.28", "4.67", "0.5342", 7),
    @("009681", "南方创新精选一年定期开放混合A", "11.24", "92.15", "3.92", "0.4406", 5),
    @("014269", "嘉实北交所精选两年定期混合A", "5.00", "87.43", "7.76", "0.3880", 1),
    @("014283", "华夏北交所创新中小企业精选两年定开混合", "3.96", "52.55", "8.26", "0.3271", 1),
    @("110012", "易方达科汇灵活配置混合", "15.73", "75.64", "2.06", "0.3240", 10),
    @("014273", "广发北交所精选两年定开混合A", "4.55", "52.69", "6.58", "0.2994", 1),
    @("014279", "汇添富北交所创新精选两年定开混合A", "4.95", "65.64", "5.62", "0.2782", 1),
    @("014294", "南方北交所精选两年定开混合", "4.63", "33.00", "4.61", "0.2134", 2),
    @("009698", "华夏成长精选6个月定期开放混合C", "3.01", "88.28", "4.67", "0.1406", 7),
    @("009682", "南方创新精选一年定期开放混合C", "3.09", "92.15", "3.92", "0.1211", 5),
    @("014274", "广发北交所精选两年定开混合C", "0.92", "52.69", "6.58", "0.0605", 1),
    @("009684", "汇添富创新增长一年定期开放混合C", "1.88", "71.21", "3.01", "0.0566", 7),
    @("014270", "嘉实北交所精选两年定期混合C", "0.64", "87.43", "7.76", "0.0497", 1),
    @("014663", "富国创新发展两年定期开放混合A", "2.62", "37.11", "1.45", "0.0380", 8),
    @("014280", "汇添富北交所创新精选两年定开混合C", "0.55", "65.64", "5.62", "0.0309", 1),
    @("014664", "富国创新发展两年定期开放混合C", "0.32", "37.11", "1.45", "0.0046", 8)
)

$r = 2
foreach ($item in $rows) {
    $ws.Range("A$r").Value = $r - 2

    $textRange = $ws.Range("B$r`:G$r")
    $textRange.NumberFormat = "@"
    $ws.Range("B$r").Value = $item[0]
    $ws.Range("C$r").Value = $item[1]
    $ws.Range("D$r").Value = $item[2]
    $ws.Range("E$r").Value = $item[3]
    $ws.Range("F$r").Value = $item[4]
    $ws.Range("G$r").Value = $item[5]

    $ws.Range("H$r").Value = $item[6]
    $r = $r + 1
}

$newTotal.Name = "总计"

# ---------------------------------------------------------------------------
# 3) Keep the originally-active sheet selected.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
